# module_pin_record.xlsx — "Added in audio, and rgb led"
#
# 1. Sheet1: re-label the old "Audio" row as "Audio buzzer", add a new
#    "RGB LED" row underneath it (each with its own "y" marker), and
#    update the dimension/selection bookkeeping that Excel recomputes
#    whenever the used range grows.
# 2. Add a new "Sheet2" after Sheet1 with a small photo/SD-card capacity
#    calculator (two input cells + two formulas), and make it the active
#    sheet/tab, matching what Excel does when a sheet is inserted while
#    the user is working in it.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: new "RGB LED" row (row 11) -----------------------------
$ws1.Range("A11").Value = "RGB LED"
$ws1.Range("V11").Value = "y"

# --- Sheet1: "Audio" (row 10) becomes "Audio buzzer", new "y" in W10 -
$ws1.Range("A10").Value = "Audio buzzer"
$ws1.Range("W10").Value = "y"

# Restore the selection to what Excel leaves it at after these edits.
$ws1.Range("V12").Select()

# --- Add Sheet2 after Sheet1 -----------------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "size of photo [kb]"
$ws2.Range("B1").Value = 50

$ws2.Range("A2").Value = "size of SD card [Gb]"
$ws2.Range("B2").Value = 4

$ws2.Range("A3").Value = "Length of time between photos [millis]"
$ws2.Range("B3").Value = 500

$ws2.Range("A4").Value = "Total # photos that can be taken"
$ws2.Range("B4").Formula = "=B2*(1000000000)/(B1*1000)"

$ws2.Range("A5").Value = "Total length of flight that can be recorded [min]"
$ws2.Range("B5").Formula = "=(B4*B3/1000)/60"

# Column A sized to fit its widest label (~44 chars wide in Excel's metric).
$ws2.Columns.Item(1).ColumnWidth = 43.1666666666667

$ws2.Range("A6").Select()
